$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 168.91667
$ws.Cells.Item(33, 9).Value = 173.2
$ws.Cells.Item(33, 10).Value = 147.5
$ws.Cells.Item(33, 11).Value = 173.2
$ws.Cells.Item(33, 12).Value = 147.5
$ws.Cells.Item(33, 13).Value = 55.80000000000001
$ws.Cells.Item(33, 14).Value = -605.5
$ws.Cells.Item(97, 8).Value = 2518.5
$ws.Cells.Item(97, 10).Value = 3405.9
$ws.Cells.Item(97, 12).Value = 10217.7
$ws.Cells.Item(97, 14).Value = -11209.7
$ws.Cells.Item(100, 8).Value = 3060.3635
$ws.Cells.Item(100, 10).Value = 3166.6667
$ws.Cells.Item(100, 12).Value = 3166.6667
$ws.Cells.Item(100, 14).Value = -4248.6667
$ws.Cells.Item(112, 8).Value = 4999.5
$ws.Cells.Item(112, 9).Value = 4999
$ws.Cells.Item(112, 10).Value = 5000
$ws.Cells.Item(112, 11).Value = 14997
$ws.Cells.Item(112, 12).Value = 15000
$ws.Cells.Item(112, 13).Value = -13889
$ws.Cells.Item(112, 14).Value = -17216
$ws.Cells.Item(137, 8).Value = 6897740
$ws.Cells.Item(137, 9).Value = 7693465.5
$ws.Cells.Item(137, 10).Value = 1450.3334
$ws.Cells.Item(137, 11).Value = 23080396.5
$ws.Cells.Item(137, 12).Value = 4351.0002
$ws.Cells.Item(137, 13).Value = -23077846.5
$ws.Cells.Item(137, 14).Value = -9451.0002
$ws.Cells.Item(138, 8).Value = 5712.7666
$ws.Cells.Item(138, 9).Value = 2627.1428
$ws.Cells.Item(138, 11).Value = 7881.428400000001
$ws.Cells.Item(138, 13).Value = -2741.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 475828.53
$ws.Cells.Item(32, 9).Value = 566380.4399999999
$ws.Cells.Item(32, 11).Value = 566380.4399999999
$ws.Cells.Item(32, 13).Value = -566093.4399999999
$ws.Cells.Item(61, 8).Value = 5510020.5
$ws.Cells.Item(61, 9).Value = 2219416
$ws.Cells.Item(61, 11).Value = 2219416
$ws.Cells.Item(61, 13).Value = -2219204
$ws.Cells.Item(74, 8).Value = 950639.8
$ws.Cells.Item(74, 9).Value = 1168526.2
$ws.Cells.Item(74, 10).Value = 6465
$ws.Cells.Item(74, 11).Value = 1168526.2
$ws.Cells.Item(74, 12).Value = 6465
$ws.Cells.Item(74, 13).Value = -1167652.2
$ws.Cells.Item(74, 14).Value = -8213
$ws.Cells.Item(77, 8).Value = 950639.8
$ws.Cells.Item(77, 9).Value = 1168526.2
$ws.Cells.Item(77, 10).Value = 6465
$ws.Cells.Item(77, 11).Value = 5842631
$ws.Cells.Item(77, 12).Value = 32325
$ws.Cells.Item(77, 13).Value = -5838263
$ws.Cells.Item(77, 14).Value = -41061
$ws.Cells.Item(102, 8).Value = 1333
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 14).Value = -5244
$ws.Cells.Item(133, 8).Value = 100000
$ws.Cells.Item(133, 9).Value = 100000
$ws.Cells.Item(133, 11).Value = 100000
$ws.Cells.Item(133, 13).Value = -97470
$ws.Cells.Item(136, 8).Value = 5510020.5
$ws.Cells.Item(136, 9).Value = 2219416
$ws.Cells.Item(136, 11).Value = 6658248
$ws.Cells.Item(136, 13).Value = -6655698

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1727.4166
$ws.Cells.Item(94, 9).Value = 1515.6207
$ws.Cells.Item(94, 10).Value = 2604.8572
$ws.Cells.Item(94, 11).Value = 1515.6207
$ws.Cells.Item(94, 12).Value = 2604.8572
$ws.Cells.Item(94, 13).Value = -1064.6207
$ws.Cells.Item(94, 14).Value = -3506.8572
$ws.Cells.Item(105, 8).Value = 1614.3334
$ws.Cells.Item(105, 9).Value = 1614.3334
$ws.Cells.Item(105, 11).Value = 1614.3334
$ws.Cells.Item(105, 13).Value = 132.6666
$ws.Cells.Item(134, 8).Value = 4833651.5
$ws.Cells.Item(134, 9).Value = 4067720
$ws.Cells.Item(134, 10).Value = 11114291
$ws.Cells.Item(134, 11).Value = 12203160
$ws.Cells.Item(134, 12).Value = 33342873
$ws.Cells.Item(134, 13).Value = -12200625
$ws.Cells.Item(134, 14).Value = -33347943

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 78669.38
$ws.Cells.Item(16, 9).Value = 1674.625
$ws.Cells.Item(16, 10).Value = 201861
$ws.Cells.Item(16, 11).Value = 1674.625
$ws.Cells.Item(16, 12).Value = 201861
$ws.Cells.Item(16, 13).Value = -1387.625
$ws.Cells.Item(16, 14).Value = -202435
$ws.Cells.Item(99, 8).Value = 12280.048
$ws.Cells.Item(99, 9).Value = 18613.75
$ws.Cells.Item(99, 11).Value = 18613.75
$ws.Cells.Item(99, 13).Value = -17115.75
$ws.Cells.Item(113, 8).Value = 78669.38
$ws.Cells.Item(113, 9).Value = 1674.625
$ws.Cells.Item(113, 10).Value = 201861
$ws.Cells.Item(113, 11).Value = 1674.625
$ws.Cells.Item(113, 12).Value = 201861
$ws.Cells.Item(113, 13).Value = 495.375
$ws.Cells.Item(113, 14).Value = -206201
$ws.Cells.Item(126, 8).Value = 12280.048
$ws.Cells.Item(126, 9).Value = 18613.75
$ws.Cells.Item(126, 11).Value = 55841.25
$ws.Cells.Item(126, 13).Value = -53371.25
$ws.Cells.Item(132, 8).Value = 3353.5715
$ws.Cells.Item(132, 9).Value = 3108.5334
$ws.Cells.Item(132, 10).Value = 3966.1667
$ws.Cells.Item(132, 11).Value = 9325.600199999999
$ws.Cells.Item(132, 12).Value = 11898.5001
$ws.Cells.Item(132, 13).Value = -6795.600199999999
$ws.Cells.Item(132, 14).Value = -16958.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1562815.8
$ws.Cells.Item(4, 9).Value = 1666976.8
$ws.Cells.Item(4, 11).Value = 5000930.4
$ws.Cells.Item(4, 13).Value = -5000818.4
$ws.Cells.Item(23, 8).Value = 250.33333
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 250.33333
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).Value = 750.99999
$ws.Cells.Item(23, 14).Value = -1220.99999
$ws.Cells.Item(107, 8).Value = 4878.7915
$ws.Cells.Item(107, 9).Value = 458.8
$ws.Cells.Item(107, 10).Value = 6041.9473
$ws.Cells.Item(107, 11).Value = 1376.4
$ws.Cells.Item(107, 12).Value = 18125.8419
$ws.Cells.Item(107, 13).Value = 543.5999999999999
$ws.Cells.Item(107, 14).Value = -21965.8419
$ws.Cells.Item(132, 8).Value = 2148
$ws.Cells.Item(132, 9).Value = 1795
$ws.Cells.Item(132, 10).Value = 2265.6667
$ws.Cells.Item(132, 11).Value = 16155
$ws.Cells.Item(132, 12).Value = 20391.0003
$ws.Cells.Item(132, 13).Value = -13625
$ws.Cells.Item(132, 14).Value = -25451.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 38799.9
$ws.Cells.Item(70, 9).Value = 66400
$ws.Cells.Item(70, 10).Value = 11199.8
$ws.Cells.Item(70, 11).Value = 66400
$ws.Cells.Item(70, 12).Value = 11199.8
$ws.Cells.Item(70, 13).Value = -66130
$ws.Cells.Item(70, 14).Value = -11739.8
$ws.Cells.Item(73, 8).Value = 38799.9
$ws.Cells.Item(73, 9).Value = 66400
$ws.Cells.Item(73, 10).Value = 11199.8
$ws.Cells.Item(73, 11).Value = 66400
$ws.Cells.Item(73, 12).Value = 11199.8
$ws.Cells.Item(73, 13).Value = -65464
$ws.Cells.Item(73, 14).Value = -13071.8
$ws.Cells.Item(122, 8).Value = 80692.234
$ws.Cells.Item(122, 9).Value = 128375.25
$ws.Cells.Item(122, 10).Value = 4399.4
$ws.Cells.Item(122, 11).Value = 385125.75
$ws.Cells.Item(122, 12).Value = 13198.2
$ws.Cells.Item(122, 13).Value = -382675.75
$ws.Cells.Item(122, 14).Value = -18098.2
$ws.Cells.Item(126, 8).Value = 2753.923
$ws.Cells.Item(126, 9).Value = 2691.75
$ws.Cells.Item(126, 11).Value = 8075.25
$ws.Cells.Item(126, 13).Value = -5605.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3559.4
$ws.Cells.Item(7, 9).Value = 3242.1428
$ws.Cells.Item(7, 10).Value = 4299.6665
$ws.Cells.Item(7, 11).Value = 3242.1428
$ws.Cells.Item(7, 12).Value = 4299.6665
$ws.Cells.Item(7, 13).Value = -3130.1428
$ws.Cells.Item(7, 14).Value = -4523.6665
$ws.Cells.Item(46, 8).Value = 3069.353
$ws.Cells.Item(46, 10).Value = 3148.6875
$ws.Cells.Item(46, 12).Value = 3148.6875
$ws.Cells.Item(46, 14).Value = -3524.6875
$ws.Cells.Item(100, 8).Value = 2866.6
$ws.Cells.Item(100, 9).Value = 1937.5
$ws.Cells.Item(100, 10).Value = 3928.4285
$ws.Cells.Item(100, 11).Value = 1937.5
$ws.Cells.Item(100, 12).Value = 3928.4285
$ws.Cells.Item(100, 13).Value = -1396.5
$ws.Cells.Item(100, 14).Value = -5010.4285
$ws.Cells.Item(122, 8).Value = 3828.125
$ws.Cells.Item(122, 9).Value = 3550
$ws.Cells.Item(122, 11).Value = 10650
$ws.Cells.Item(122, 13).Value = -8200
$ws.Cells.Item(126, 8).Value = 3559.4
$ws.Cells.Item(126, 9).Value = 3242.1428
$ws.Cells.Item(126, 10).Value = 4299.6665
$ws.Cells.Item(126, 11).Value = 9726.428400000001
$ws.Cells.Item(126, 12).Value = 12898.9995
$ws.Cells.Item(126, 13).Value = -7256.428400000001
$ws.Cells.Item(126, 14).Value = -17838.9995
$ws.Cells.Item(132, 8).Value = 761091.5
$ws.Cells.Item(132, 9).Value = 1152864.5
$ws.Cells.Item(132, 10).Value = 3663.6667
$ws.Cells.Item(132, 11).Value = 3458593.5
$ws.Cells.Item(132, 12).Value = 10991.0001
$ws.Cells.Item(132, 13).Value = -3456063.5
$ws.Cells.Item(132, 14).Value = -16051.0001
$ws.Cells.Item(136, 8).Value = 12154519
$ws.Cells.Item(136, 9).Value = 10715562
$ws.Cells.Item(136, 10).Value = 22227220
$ws.Cells.Item(136, 11).Value = 32146686
$ws.Cells.Item(136, 12).Value = 66681660
$ws.Cells.Item(136, 13).Value = -32144136
$ws.Cells.Item(136, 14).Value = -66686760

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3764
$ws.Cells.Item(81, 10).Value = 3824.0833
$ws.Cells.Item(81, 12).Value = 7648.1666
$ws.Cells.Item(81, 14).Value = -9770.1666
$ws.Cells.Item(84, 8).Value = 3764
$ws.Cells.Item(84, 10).Value = 3824.0833
$ws.Cells.Item(84, 12).Value = 38240.833
$ws.Cells.Item(84, 14).Value = -48848.833
$ws.Cells.Item(113, 8).Value = 1705.5416
$ws.Cells.Item(113, 9).Value = 1756.75
$ws.Cells.Item(113, 11).Value = 5270.25
$ws.Cells.Item(113, 13).Value = -3100.25
$ws.Cells.Item(126, 8).Value = 2343.4167
$ws.Cells.Item(126, 9).Value = 2323
$ws.Cells.Item(126, 10).Value = 2445.5
$ws.Cells.Item(126, 11).Value = 6969
$ws.Cells.Item(126, 12).Value = 7336.5
$ws.Cells.Item(126, 13).Value = -4499
$ws.Cells.Item(126, 14).Value = -12276.5
$ws.Cells.Item(132, 8).Value = 3970892.2
$ws.Cells.Item(132, 9).Value = 4388528.5
$ws.Cells.Item(132, 10).Value = 3349.75
$ws.Cells.Item(132, 11).Value = 13165585.5
$ws.Cells.Item(132, 12).Value = 10049.25
$ws.Cells.Item(132, 13).Value = -13163055.5
$ws.Cells.Item(132, 14).Value = -15109.25
$ws.Cells.Item(136, 8).Value = 1768567.2
$ws.Cells.Item(136, 9).Value = 1489689.9
$ws.Cells.Item(136, 11).Value = 4469069.699999999
$ws.Cells.Item(136, 13).Value = -4469069.699999999
